# REPORTGEN-455 : New component for Top Components by properties
#
# Adds a new worksheet "3 - Top Components" at the end of the workbook,
# documenting the new TOP_COMPONENTS_BY_PROPERTIES placeholder, following
# the same layout/style conventions as the other "3 - ..." documentation
# sheets (e.g. "3 - Evolution of standards").

$wb = $excel.ActiveWorkbook

# Use the existing "Evolution of standards" sheet as the style template -
# it has the exact title / body / notes formatting we want to reuse.
$template = $wb.Worksheets.Item("3 - Evolution of standards")

# Insert the new sheet right after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "3 - Top Components"

# --- Column / header layout -------------------------------------------------
$ws.Columns(2).ColumnWidth = 121.29

# --- Content (shared strings are appended in this exact order) -------------
$ws.Range("B1").Value = "3.9.`tTop components by properties"
$ws.Range("B2").Value = "* Block Name = TOP_COMPONENTS_BY_PROPERTIES"
$ws.Range("B4").Value = "PROP1 : name of first property, cyclomaticComplexity if not exists"
$ws.Range("B5").Value = "PROP2 : name of second property, fanOut if not exists"
$ws.Range("B6").Value = "ORDER1 : ASC or DESC for PROP1, DESC by default"
$ws.Range("B7").Value = "ORDER2 : ASC or DESC for PROP2, DESC by default"
$ws.Range("B8").Value = "COUNT: the number of lines to display, 50 by default (-1 or all is not allowed, it will take too much time and paper)"
$ws.Range("B3").Value = "* Options :"
$ws.Range("B10").Value = "For PROP1 and PROP2, the available values are : codeLines, commentedCodeLines, commentLines, coupling, fanIn, fanOut, cyclomaticComplexity, ratioCommentLinesCodeLines, halsteadProgramLength, halsteadProgramVocabulary, halsteadVolume, distinctOperators, distinctOperands, integrationComplexity, essentialComplexity"
$ws.Range("B15").Value = "RepGen:TABLE;TOP_COMPONENTS_BY_PROPERTIES;PROP1=cyclomaticComplexity,PROP2=ratioCommentLinesCodeLines,ORDER1=desc,ORDER2=asc,COUNT=10"
$ws.Range("B11").Value = "If PROP1 and/or PROP2 is not correctly set,list of available values is displayed"
$ws.Range("B13").Value = "Note : This component is only relevant on an engineering database. It is empty on an analytics database."

# --- Formatting: reuse the styles already used on the template sheet -------
# Title (row 1, bold 14pt)
$template.Range("B1").Copy()
$ws.Range("B1").PasteSpecial(-4122)

# Body text (rows 2-4, regular 11pt)
$template.Range("B2").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)

# Option list (rows 5-8, regular 12pt)
$template.Range("B4").Copy()
$ws.Range("B5:B8").PasteSpecial(-4122)

# Blank spacer row (row 9, regular 12pt - same family as the option list)
$template.Range("B8").Copy()
$ws.Range("B9").PasteSpecial(-4122)

$ws.Application.CutCopyMode = $false

# Row heights to match the template sheet
$ws.Rows(1).RowHeight = 18.75
$ws.Rows(5).RowHeight = 15.75
$ws.Rows(6).RowHeight = 15.75
$ws.Rows(7).RowHeight = 15.75
$ws.Rows(8).RowHeight = 15.75
$ws.Rows(9).RowHeight = 15.75

# --- Selection / activation -------------------------------------------------
$ws.Range("B11").Select()
$ws.Activate()

Write-Host "Added sheet '3 - Top Components'"
